$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("Material cost")
$ws5.Columns("B").ColumnWidth = 11.1328125
